# Deploy the implementation guide:
#  - refresh the generation Date
#  - refresh the Contact display text
#  - add a new "Jurisdiction" metadata row (inserted right after "Contact")
#
# The "Concepts" sheet (sheet2) needs no direct edits: its content is
# unaffected by this change (only the underlying shared-string table grows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update the "Date" row (row 8) -----------------------------------------
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# --- Update the "Contact" row (row 10) --------------------------------------
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# --- Insert a new "Jurisdiction" row right after "Contact" (row 11) --------
# This pushes "Description" and every following row down by one.
$ws.Rows.Item(11).Insert()

# Re-apply the same formatting used by the rest of the data rows (borders,
# top-aligned wrapped text) so the new row matches its neighbours instead of
# picking up a blank/default style.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
